$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "1047414864"
$ws.Range("D16").Value = "KARINA PAOLA BLANCO GAMARRA"
$ws.Range("E16").Value = "2306"
$ws.Range("F16").Value = 46400

$ws.Range("C17").Value = "33102174"
$ws.Range("D17").Value = "ANA MARINA CALVO CARTAGENA"
$ws.Range("E17").Value = "2306"
$ws.Range("F17").Value = 46400

$ws.Range("C18").Value = "1052951335"
$ws.Range("D18").Value = "AMALIA CANDELARIA GARCIA SANDOVAL"
$ws.Range("E18").Value = "2306"
$ws.Range("F18").Value = 46400

$ws.Range("C19").Value = "1143377977"
$ws.Range("D19").Value = "JULIETH NILET BARRIOS MAUSSA"
$ws.Range("E19").Value = "2306"
$ws.Range("F19").Value = 46400

$ws.Range("C20").Value = "1128063714"
$ws.Range("D20").Value = "KATHERINE RUIZ BLANCO"
$ws.Range("E20").Value = "2306"
$ws.Range("F20").Value = 46400

$ws.Range("C21").Value = "1050950839"
$ws.Range("D21").Value = "SINDY PAOLA HERNANDEZ MACHACON"
$ws.Range("E21").Value = "2306"
$ws.Range("F21").Value = 46400

$ws.Range("C22").Value = "1143325274"
$ws.Range("D22").Value = "XIOMARA SALGADO CASSIANI"
$ws.Range("E22").Value = "2306"
$ws.Range("F22").Value = 46400

$ws.Range("C23").Value = "1049942585"
$ws.Range("D23").Value = "YURIS JOHANA CUTEN JULIO"
$ws.Range("E23").Value = "2306"
$ws.Range("F23").Value = 46400

$ws.Range("C24").Value = "1102830646"
$ws.Range("D24").Value = "GISSETH DEL CARMEN DECHAMP MORALES"
$ws.Range("E24").Value = "2312"
$ws.Range("F24").Value = 25333

$ws.Range("C25").Value = "1102830646"
$ws.Range("D25").Value = "GISSETH DEL CARMEN DECHAMP MORALES"
$ws.Range("E25").Value = "2311"
$ws.Range("F25").Value = 40000

$ws.Range("C26").Value = "1102830646"
$ws.Range("D26").Value = "GISSETH DEL CARMEN DECHAMP MORALES"
$ws.Range("E26").Value = "2310"
$ws.Range("F26").Value = 46400

$ws.Range("C27").Value = "1102830646"
$ws.Range("D27").Value = "GISSETH DEL CARMEN DECHAMP MORALES"
$ws.Range("E27").Value = "2309"
$ws.Range("F27").Value = 46400

$ws.Range("C28").Value = "1102830646"
$ws.Range("D28").Value = "GISSETH DEL CARMEN DECHAMP MORALES"
$ws.Range("E28").Value = "2308"
$ws.Range("F28").Value = 46400

$ws.Range("C29").Value = "1102830646"
$ws.Range("D29").Value = "GISSETH DEL CARMEN DECHAMP MORALES"
$ws.Range("E29").Value = "2307"
$ws.Range("F29").Value = 46400

$ws.Range("C30").Value = "1102830646"
$ws.Range("D30").Value = "GISSETH DEL CARMEN DECHAMP MORALES"
$ws.Range("E30").Value = "2306"
$ws.Range("F30").Value = 46400

$ws.Range("C31").Value = "45564420"
$ws.Range("D31").Value = "SANDY HELENA PUPO LEON"
$ws.Range("E31").Value = "2306"
$ws.Range("F31").Value = 46400

$ws.Range("C32").Value = "1047496345"
$ws.Range("D32").Value = "KEYRA LUZ NEGRETTE BAUTISTA"
$ws.Range("E32").Value = "2306"
$ws.Range("F32").Value = 46400
